$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Add the three new character styles used by the new runs.
# ---------------------------------------------------------------------

# wdStyleTypeCharacter = 2
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------
# 2. Apply GaNParagraph to every run containing the recurring Swedish
#    "Du deltar..." paragraph (5 occurrences in the document).
# ---------------------------------------------------------------------

$paragraphText = "Du deltar i en världsomspännande kampanj*den synliga natthimlens över hela världen."

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($paragraphText, $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNParagraph"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($paragraphText, $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
}

# ---------------------------------------------------------------------
# 3. Apply GaNLinks to the run with the GaNight map link.
# ---------------------------------------------------------------------

$linkRng = $d.Content
$linkRng.Find.ClearFormatting()
$linkFound = $linkRng.Find.Execute("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($linkFound) {
    $linkRng.Style = "GaNLinks"
}
